$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 and C5 swap content: the "Triangle (Transform the mockup database...)" task
# moves from C2 down to C5, and the "Circle (Continue with analysis...)" task
# moves from C5 up to C2 (also fixing a stray double-space typo in the Triangle text).
$ws.Range("C2").Value = "Circle (Continue with analysis and `ncreate visuals to accompany the data story.)"
$ws.Range("C5").Value = "Triangle (Transform the mockup database into a full database that integrates with your work.)"

# Both cells pick up wrap-text plus an explicitly-applied (visually default)
# font, mirroring the new cellXfs entry added for this edit.
$ws.Range("C2").WrapText = $true
$ws.Range("C2").Font.ThemeColor = 1
$ws.Range("C5").WrapText = $true
$ws.Range("C5").Font.ThemeColor = 1

# View / column tweaks
$ws.Columns("E").ColumnWidth = 54.6640625
$excel.ActiveWindow.Zoom = 144
$ws.Range("C2").Select()
